# Updating latest mobile code
# Update the "wild1_instance" row (row 14):
#  - H14 keeps referencing the same logical string, but its text content is updated
#    to include the additional ranges "20-21,24-31,35-".
#  - G14 gets a brand new value, a different range string that did not exist before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update H14 text (existing shared string's content changes)
$ws.Range("H14").Value = "1-10,12-16,20-21,24-31,35-37,39-43,46-54,56-60,62-63"

# Update G14 to a brand new string value
$ws.Range("G14").Value = "1-10,12-14,20-21,24-31,37,39-43,46-54,56-60,62-63"
